$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("choices")
$ws3 = $wb.Worksheets.Item("settings")

# --- Fix B18/B19 format (lose the "horizontal right" alignment, match rest of column) ---
$ws2.Range("B17").Copy()
$ws2.Range("B18").PasteSpecial(-4122)
$ws2.Range("B17").Copy()
$ws2.Range("B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fix B1 header format (lose vertical-center, keep wrap only) ---
$ws2.Range("B2").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Convert all choice values in column B from numbers to "a"-prefixed text ---
$ws2.Range("B2").Value = "a1"
$ws2.Range("B3").Value = "a0"
$ws2.Range("B4").Value = "a888"
$ws2.Range("B5").Value = "a1"
$ws2.Range("B6").Value = "a0"
$ws2.Range("B7").Value = "a888"
$ws2.Range("B8").Value = "a888"
$ws2.Range("B9").Value = "a999"
$ws2.Range("B10").Value = "a1"
$ws2.Range("B11").Value = "a2"
$ws2.Range("B12").Value = "a1"
$ws2.Range("B13").Value = "a2"
$ws2.Range("B14").Value = "a3"
$ws2.Range("B15").Value = "a4"
$ws2.Range("B16").Value = "a5"
$ws2.Range("B17").Value = "a6"
$ws2.Range("B18").Value = "a7"
$ws2.Range("B19").Value = "a888"
$ws2.Range("B20").Value = "a1"
$ws2.Range("B21").Value = "a2"
$ws2.Range("B22").Value = "a3"
$ws2.Range("B23").Value = "a4"
$ws2.Range("B24").Value = "a5"
$ws2.Range("B25").Value = "a6"
$ws2.Range("B26").Value = "a7"
$ws2.Range("B27").Value = "a8"
$ws2.Range("B28").Value = "a9"
$ws2.Range("B29").Value = "a10"
$ws2.Range("B30").Value = "a888"
$ws2.Range("B31").Value = "a0"
$ws2.Range("B32").Value = "a1"
$ws2.Range("B33").Value = "a2"
$ws2.Range("B34").Value = "a3"
$ws2.Range("B35").Value = "a4"
$ws2.Range("B36").Value = "a5"
$ws2.Range("B37").Value = "a6"
$ws2.Range("B38").Value = "a7"
$ws2.Range("B39").Value = "a888"
$ws2.Range("B40").Value = "a0"
$ws2.Range("B41").Value = "a1"
$ws2.Range("B42").Value = "a2"
$ws2.Range("B43").Value = "a3"
$ws2.Range("B44").Value = "a4"
$ws2.Range("B45").Value = "a888"
$ws2.Range("B46").Value = "a0"
$ws2.Range("B47").Value = "a1"
$ws2.Range("B48").Value = "a2"
$ws2.Range("B49").Value = "a3"
$ws2.Range("B50").Value = "a4"
$ws2.Range("B51").Value = "a5"
$ws2.Range("B52").Value = "a888"
$ws2.Range("B53").Value = "a1"
$ws2.Range("B54").Value = "a2"
$ws2.Range("B55").Value = "a3"
$ws2.Range("B56").Value = "a4"
$ws2.Range("B57").Value = "a1"
$ws2.Range("B58").Value = "a2"
$ws2.Range("B59").Value = "a3"
$ws2.Range("B60").Value = "a999"
$ws2.Range("B61").Value = "a0"
$ws2.Range("B62").Value = "a1"
$ws2.Range("B63").Value = "a888"
$ws2.Range("B64").Value = "a999"
$ws2.Range("B65").Value = "a1"
$ws2.Range("B66").Value = "a0"
$ws2.Range("B67").Value = "a888"
$ws2.Range("B68").Value = "a999"
$ws2.Range("B69").Value = "a-777"
$ws2.Range("B70").Value = "a-888"
$ws2.Range("B71").Value = "a-999"
$ws2.Range("B72").Value = "a1"
$ws2.Range("B73").Value = "a0"
$ws2.Range("B74").Value = "a2"
$ws2.Range("B75").Value = "a999"
$ws2.Range("B76").Value = "a1"
$ws2.Range("B77").Value = "a2"
$ws2.Range("B78").Value = "a3"
$ws2.Range("B79").Value = "a1"
$ws2.Range("B80").Value = "a2"
$ws2.Range("B81").Value = "a3"
$ws2.Range("B82").Value = "a4"
$ws2.Range("B83").Value = "a5"
$ws2.Range("B84").Value = "a6"
$ws2.Range("B85").Value = "a888"
$ws2.Range("B86").Value = "a0"
$ws2.Range("B87").Value = "a1"
$ws2.Range("B88").Value = "a1"
$ws2.Range("B89").Value = "a1"
$ws2.Range("B90").Value = "a0"
$ws2.Range("B91").Value = "a1"
$ws2.Range("B92").Value = "a0"
$ws2.Range("B93").Value = "a1"
$ws2.Range("B94").Value = "a0"
$ws2.Range("B95").Value = "a2"

# --- Column B width: reset to default sheet width (was custom 18.66) ---
$ws2.Columns.Item(2).ColumnWidth = 9.997968749999998

# --- Active sheet / selection: choices becomes the active tab ---
$ws2.Select()
$ws2.Range("E6").Select()
